$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Index (A), Time (B), SegStart (C), SegEnd (D), AWC (E)
$data = @(
    @(2, 39, "5:45 AM", 1169.48, 1199.48, 21.87),
    @(3, 91, "6:11 AM", 2725.88, 2755.88, 22.16),
    @(4, 100, "6:15 AM", 2986.3175, 3016.3175, 30.9975),
    @(5, 120, "6:25 AM", 3577.16, 3607.16, 49.37),
    @(6, 151, "6:40 AM", 4506.15, 4536.15, 33.81),
    @(7, 174, "6:52 AM", 5193.04, 5223.04, 71.97),
    @(8, 182, "6:56 AM", 5440.89, 5470.89, 108.74),
    @(9, 230, "7:20 AM", 6880.786667, 6910.786667, 32.193333),
    @(10, 239, "7:24 AM", 7146.62, 7176.62, 32.95),
    @(11, 251, "7:30 AM", 7516.66, 7546.66, 27.63),
    @(12, 275, "7:43 AM", 8248.530000000001, 8278.530000000001, 21.35),
    @(13, 331, "8:10 AM", 9903.9, 9933.9, 66.37),
    @(14, 524, "9:47 AM", 15702.116667, 15732.116667, 35.09),
    @(15, 550, "10:00 A", 16484, 16514, 26.08),
    @(16, 570, "10:10 A", 17087.03, 17117.03, 37.05),
    @(17, 993, "1:41 PM", 29775.52, 29805.52, 41.35),
    @(18, 1016, "1:53 PM", 30454.73, 30484.73, 24.88),
    @(19, 1022, "1:56 PM", 30644.296667, 30674.296667, 24.493333),
    @(20, 1036, "2:03 PM", 31058.49, 31088.49, 60.24),
    @(21, 1086, "2:28 PM", 32575.77, 32605.77, 21.33),
    @(22, 1103, "2:36 PM", 33070.88, 33100.88, 25.1),
    @(23, 1141, "2:55 PM", 34219.87, 34249.87, 21.48),
    @(24, 1185, "3:17 PM", 35530.52, 35560.52, 115.925),
    @(25, 1196, "3:23 PM", 35858.715, 35888.715, 56.205),
    @(26, 1202, "3:26 PM", 36041.01, 36071.01, 24.126667),
    @(27, 1209, "3:29 PM", 36244.02, 36274.02, 75.29000000000001),
    @(28, 1221, "3:35 PM", 36616.4, 36646.4, 23.47),
    @(29, 1227, "3:39 PM", 36808.87, 36838.87, 22.12),
    @(30, 1245, "3:47 PM", 37332.09, 37362.09, 193.12),
    @(31, 1257, "3:53 PM", 37684.8, 37714.8, 36.82),
    @(32, 1273, "4:02 PM", 38186.04, 38216.04, 28.77),
    @(33, 1280, "4:05 PM", 38387.28, 38417.28, 28.405),
    @(34, 1291, "4:11 PM", 38728.38, 38758.38, 100.57),
    @(35, 1300, "4:15 PM", 38978.615, 39008.615, 27.355),
    @(36, 1314, "4:22 PM", 39391.77, 39421.77, 27.71),
    @(37, 1333, "4:31 PM", 39975.185, 40005.185, 26.385),
    @(38, 1343, "4:36 PM", 40265.54, 40295.54, 30.27),
    @(39, 1362, "4:46 PM", 40838.935, 40868.935, 104.33),
    @(40, 1369, "4:49 PM", 41046.07, 41076.07, 33.66),
    @(41, 1390, "5:00 PM", 41680.48, 41710.48, 38.76),
    @(42, 1412, "5:11 PM", 42334.07, 42364.07, 101.66),
    @(43, 1418, "5:14 PM", 42521.595, 42551.595, 25.965),
    @(44, 1427, "5:19 PM", 42809.97, 42839.97, 161.82),
    @(45, 1445, "5:27 PM", 43323.345, 43353.345, 29.31),
    @(46, 1451, "5:30 PM", 43516.02, 43546.02, 88.58),
    @(47, 1472, "5:41 PM", 44138.865, 44168.865, 48.69),
    @(48, 1479, "5:44 PM", 44340.68, 44370.68, 49.65),
    @(49, 1492, "5:51 PM", 44738.87, 44768.87, 22.62),
    @(50, 1522, "6:06 PM", 45637.495, 45667.495, 33.28),
    @(51, 1551, "6:20 PM", 46502.28, 46532.28, 39.27)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Output "Updated $($data.Count) rows"
